$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that *looks* numeric (e.g. "389.28") as literal text,
# matching the sheet's existing inline-string "Price" column formatting,
# without leaving a residual custom number format behind on the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '51.733.95'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '3.093.99'
$ws.Range("E3").Value = '  +3.74%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws.Range("D5") '389.28'
$ws.Range("E5").Value = '  +2.00%  '
Set-TextValue $ws.Range("D6") '103.63'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.02%  '
Set-TextValue $ws.Range("D9") '0.588'
$ws.Range("E9").Value = '  -0.56%  '
Set-TextValue $ws.Range("D10") '37.09'
$ws.Range("E10").Value = '  +1.27%  '
Set-TextValue $ws.Range("D12") '0.0861'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").Value = '3.584.61'
$ws.Range("E13").Value = '  +3.71%  '
Set-TextValue $ws.Range("D14") '18.70'
$ws.Range("E14").Value = '  +1.43%  '
Set-TextValue $ws.Range("D15") '7.81'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '3.093.69'
$ws.Range("E16").Value = '  +3.71%  '
Set-TextValue $ws.Range("D17") '0.984'
$ws.Range("E17").Value = '  -1.24%  '
Set-TextValue $ws.Range("D18") '10.64'
$ws.Range("E18").Value = '  -4.62%  '
$ws.Range("D19").Value = '51.871.83'
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("E20").Value = '  +2.28%  '
Set-TextValue $ws.Range("D21") '12.50'
$ws.Range("E21").Value = '  -0.68%  '
Set-TextValue $ws.Range("D23") '70.04'
$ws.Range("E23").Value = '  -0.29%  '
Set-TextValue $ws.Range("D24") '268.79'
$ws.Range("E25").Value = '  -2.45%  '
Set-TextValue $ws.Range("D26") '8.19'
$ws.Range("E26").Value = '  +4.72%  '
Set-TextValue $ws.Range("D27") '27.15'
$ws.Range("E27").Value = '  +4.21%  '
$ws.Range("E28").Value = '  +1.37%  '
Set-TextValue $ws.Range("D29") '7.25'
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("E32").Value = '  -0.13%  '
Set-TextValue $ws.Range("D33") '35.75'
$ws.Range("E33").Value = '  +2.96%  '
Set-TextValue $ws.Range("D34") '2.06'
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("E37").Value = '  -0.13%  '
Set-TextValue $ws.Range("D38") '3.40'
$ws.Range("E38").Value = '  +3.39%  '
Set-TextValue $ws.Range("D39") '0.294'
$ws.Range("E39").Value = '  +8.43%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D40") '17.05'
$ws.Range("E40").Value = '  +0.68%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D41") '1.89'
$ws.Range("E41").Value = '  +2.19%  '
$ws.Range("E42").Value = '  +0.14%  '
Set-TextValue $ws.Range("D44") '127.33'
$ws.Range("E44").Value = '  +1.40%  '
$ws.Range("E45").Value = '  -3.30%  '
Set-TextValue $ws.Range("D46") '22.21'
$ws.Range("E46").Value = '  +3.79%  '
Set-TextValue $ws.Range("D47") '2.47'
$ws.Range("E47").Value = '  +4.24%  '
$ws.Range("E48").Value = '  +2.25%  '
$ws.Range("D49").Value = '2.048.72'
$ws.Range("E49").Value = '  +1.24%  '
$ws.Range("D50").Value = '3.403.27'
$ws.Range("E50").Value = '  +3.75%  '
Set-TextValue $ws.Range("D51") '0.207'
$ws.Range("E51").Value = '  +6.73%  '